# Auto-generated Excel COM-interop script applying market-data refresh
# updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2222
$ws.Range("J40").Value = 2222
$ws.Range("L40").Value = 2222
$ws.Range("N40").Value = -2572

$ws.Range("H64").Value = 6149
$ws.Range("I64").Value = 6579
$ws.Range("K64").Value = 6579
$ws.Range("M64").Value = -6331

$ws.Range("H67").Value = 6149
$ws.Range("I67").Value = 6579
$ws.Range("K67").Value = 6579
$ws.Range("M67").Value = -5721

$ws.Range("H74").Value = 19237076
$ws.Range("I74").Value = 6549
$ws.Range("K74").Value = 6549
$ws.Range("M74").Value = -5613

$ws.Range("H77").Value = 19237076
$ws.Range("I77").Value = 6549
$ws.Range("K77").Value = 32745
$ws.Range("M77").Value = -28065

$ws.Range("H132").Value = 1902.7715
$ws.Range("I132").Value = 1635.08
$ws.Range("K132").Value = 4905.24
$ws.Range("M132").Value = -2375.24

$ws.Range("H138").Value = 7163.3
$ws.Range("I138").Value = 12947.223
$ws.Range("J138").Value = 4684.476
$ws.Range("K138").Value = 38841.669
$ws.Range("L138").Value = 14053.428
$ws.Range("M138").Value = -33701.669
$ws.Range("N138").Value = -24333.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 300980.28
$ws.Range("I74").Value = 2155.0164
$ws.Range("J74").Value = 1373235.6
$ws.Range("K74").Value = 2155.0164
$ws.Range("L74").Value = 1373235.6
$ws.Range("M74").Value = -1281.0164
$ws.Range("N74").Value = -1374983.6

$ws.Range("H77").Value = 300980.28
$ws.Range("I77").Value = 2155.0164
$ws.Range("J77").Value = 1373235.6
$ws.Range("K77").Value = 10775.082
$ws.Range("L77").Value = 6866178
$ws.Range("M77").Value = -6407.082
$ws.Range("N77").Value = -6874914

$ws.Range("H122").Value = 3202.197
$ws.Range("I122").Value = 2220.7354
$ws.Range("K122").Value = 6662.206200000001
$ws.Range("M122").Value = -4212.206200000001

$ws.Range("H132").Value = 1855.1515
$ws.Range("I132").Value = 1473.1034
$ws.Range("K132").Value = 4419.3102
$ws.Range("M132").Value = -1889.3102

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1050.5555
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5346

$ws.Range("H64").Value = 1428.6666
$ws.Range("I64").Value = 1238.3334
$ws.Range("J64").Value = 1492.1111
$ws.Range("K64").Value = 1238.3334
$ws.Range("L64").Value = 1492.1111
$ws.Range("M64").Value = -1013.3334
$ws.Range("N64").Value = -1942.1111

$ws.Range("H67").Value = 1428.6666
$ws.Range("I67").Value = 1238.3334
$ws.Range("J67").Value = 1492.1111
$ws.Range("K67").Value = 1238.3334
$ws.Range("L67").Value = 1492.1111
$ws.Range("M67").Value = -458.3334
$ws.Range("N67").Value = -3052.1111

$ws.Range("H86").Value = 5152.3687
$ws.Range("I86").Value = 3755.8333
$ws.Range("J86").Value = 7546.4287
$ws.Range("K86").Value = 3755.8333
$ws.Range("L86").Value = 7546.4287
$ws.Range("M86").Value = -2632.8333
$ws.Range("N86").Value = -9792.4287

$ws.Range("H89").Value = 5152.3687
$ws.Range("I89").Value = 3755.8333
$ws.Range("J89").Value = 7546.4287
$ws.Range("K89").Value = 18779.1665
$ws.Range("L89").Value = 37732.14350000001
$ws.Range("M89").Value = -13163.1665
$ws.Range("N89").Value = -48964.14350000001

$ws.Range("H134").Value = 23686506
$ws.Range("I134").Value = 2111.8064
$ws.Range("J134").Value = 128574536
$ws.Range("K134").Value = 6335.4192
$ws.Range("L134").Value = 385723608
$ws.Range("M134").Value = -3800.4192
$ws.Range("N134").Value = -385728678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1279.75
$ws.Range("I22").Value = 1279.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1279.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -929.75
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 2364.528
$ws.Range("I31").Value = 1689.25
$ws.Range("K31").Value = 1689.25
$ws.Range("M31").Value = -1394.25

$ws.Range("H34").Value = 2364.528
$ws.Range("I34").Value = 1689.25
$ws.Range("K34").Value = 1689.25
$ws.Range("M34").Value = -1487.25

$ws.Range("H109").Value = 45500
$ws.Range("J109").Value = 45500
$ws.Range("L109").Value = 45500
$ws.Range("N109").Value = -47580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 9067
$ws.Range("J54").Value = 11600.5
$ws.Range("L54").Value = 34801.5
$ws.Range("N54").Value = -35919.5

$ws.Range("H68").Value = 1707.1538
$ws.Range("J68").Value = 1704.5902
$ws.Range("L68").Value = 5113.7706
$ws.Range("N68").Value = -6735.7706

$ws.Range("H71").Value = 1707.1538
$ws.Range("J71").Value = 1704.5902
$ws.Range("L71").Value = 15341.3118
$ws.Range("N71").Value = -23453.3118

$ws.Range("H95").Value = 18997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 63.357143
$ws.Range("I2").Value = 63.357143
$ws.Range("K2").Value = 63.357143
$ws.Range("M2").Value = 49.642857

$ws.Range("H97").Value = 716.6667
$ws.Range("I97").Value = 660
$ws.Range("K97").Value = 660
$ws.Range("M97").Value = -164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12925.363
$ws.Range("J7").Value = 14622.625
$ws.Range("L7").Value = 14622.625
$ws.Range("N7").Value = -14846.625

$ws.Range("H55").Value = 500461.12
$ws.Range("I55").Value = 1667133.1
$ws.Range("K55").Value = 1667133.1
$ws.Range("M55").Value = -1666960.1

$ws.Range("H61").Value = 3898.25
$ws.Range("I61").Value = 1598
$ws.Range("K61").Value = 1598
$ws.Range("M61").Value = -1396

$ws.Range("H93").Value = 1390.1333
$ws.Range("I93").Value = 1511.5385
$ws.Range("K93").Value = 1511.5385
$ws.Range("M93").Value = -263.5385000000001

$ws.Range("H98").Value = 43333.332
$ws.Range("J98").Value = 43333.332
$ws.Range("L98").Value = 43333.332
$ws.Range("N98").Value = -49323.332

$ws.Range("H113").Value = 3898.25
$ws.Range("I113").Value = 1598
$ws.Range("K113").Value = 1598
$ws.Range("M113").Value = 572

$ws.Range("H126").Value = 12925.363
$ws.Range("J126").Value = 14622.625
$ws.Range("L126").Value = 43867.875
$ws.Range("N126").Value = -48807.875

$ws.Range("H132").Value = 3586.9
$ws.Range("J132").Value = 6496.2856
$ws.Range("L132").Value = 19488.8568
$ws.Range("N132").Value = -24548.8568

$ws.Range("H136").Value = 2059.845
$ws.Range("J136").Value = 2582.4666
$ws.Range("L136").Value = 7747.399800000001
$ws.Range("N136").Value = -12847.3998

$ws.Range("H141").Value = 137310
$ws.Range("J141").Value = 137310
$ws.Range("L141").Value = 137310
$ws.Range("N141").Value = -147670

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H122").Value = 2505.4285
$ws.Range("I122").Value = 2505.4285
$ws.Range("K122").Value = 7516.2855
$ws.Range("M122").Value = -5066.2855

$ws.Range("H126").Value = 1880.3572
$ws.Range("I126").Value = 2049
$ws.Range("K126").Value = 6147
$ws.Range("M126").Value = -3677

$ws.Range("H132").Value = 2102.6538
$ws.Range("I132").Value = 1804.5294
$ws.Range("K132").Value = 5413.5882
$ws.Range("M132").Value = -2883.5882

$ws.Range("H133").Value = 79000
$ws.Range("J133").Value = 79000
$ws.Range("L133").Value = 79000
$ws.Range("N133").Value = -89120

$ws.Range("H135").Value = 88165.336
$ws.Range("J135").Value = 88165.336
$ws.Range("L135").Value = 88165.336
$ws.Range("N135").Value = -98305.336
